$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the truncated two-digit years in the START_DATE / END_DATE columns
# to full four-digit years (enables filtering studies by year).
$ws.Range("D2").Value = "2017/01/17"
$ws.Range("E2").Value = "2017/12/01"
$ws.Range("D3").Value = "2019/01/01"
$ws.Range("E3").Value = "2019/12/01"
$ws.Range("D4").Value = "2017/01/17"
$ws.Range("E4").Value = "2017/12/01"
$ws.Range("D5").Value = "2019/01/01"
$ws.Range("E5").Value = "2019/12/01"

# Move the active selection to E5, matching the saved view state.
$ws.Range("E5").Select()
